$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.093.07"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "1.665.67"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.59"
$ws.Range("E5").Value = "  -3.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5248"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("E7").Value = "  -0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2623"
$ws.Range("E8").Value = "  -3.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06286"
$ws.Range("E9").Value = "  -1.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.12"
$ws.Range("E10").Value = "  -2.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07536"
$ws.Range("E11").Value = "  -1.75%  "

$ws.Range("D12").Value = "1.666.03"
$ws.Range("E12").Value = "  -1.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.435"
$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5523"
$ws.Range("E14").Value = "  -4.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.49"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000007910"
$ws.Range("E16").Value = "  -4.93%  "

$ws.Range("D17").Value = "26.131.10"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.720"
$ws.Range("E19").Value = "  -3.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.33"
$ws.Range("E20").Value = "  -2.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("E21").Value = "  -4.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.153"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.62"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1245"
$ws.Range("E25").Value = "  -2.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.471"
$ws.Range("E26").Value = "  -4.51%  "

$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06306"
$ws.Range("E28").Value = "  +2.42%  "

$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.272"
$ws.Range("E30").Value = "  -3.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.487"
$ws.Range("E31").Value = "  -2.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.404"
$ws.Range("E32").Value = "  -4.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.632"
$ws.Range("E33").Value = "  -2.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9952"
$ws.Range("E34").Value = "  -3.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6017"
$ws.Range("E35").Value = "  -3.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.406"
$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").Value = "1.106.40"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.076"
$ws.Range("E39").Value = "  -0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01611"
$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8697"
$ws.Range("E41").Value = "  -0.99%  "

$ws.Range("E42").Value = "  -0.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.72"
$ws.Range("E43").Value = "  -1.04%  "

$ws.Range("D44").Value = "1.818.50"
$ws.Range("E44").Value = "  -1.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -2.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.31"
$ws.Range("E46").Value = "  -3.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.002"
$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05235"
$ws.Range("E49").Value = "  -0.89%  "

$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.946"
$ws.Range("E51").Value = "  -1.73%  "
